# The sheet is being re-shaped: two new columns ("ownTeam" / "oppTeam") are
# inserted before "batsman", a new match row ("Dubai (DSC)") is added, and
# the remaining rows are re-ordered chronologically. Simplest + most robust
# way to land the new, final state is to clear the used range and rewrite
# every cell explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used after "Isuru Udana" in the source data.
$nbsp = [char]0x00A0
$batsman = "Isuru Udana" + $nbsp

# Wipe whatever is currently on the sheet before laying out the new table.
$ws.Cells.Clear()

$headers = @("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")

$data = @(
    @(" Sharjah", " October 15 2020", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", $batsman, "10", "5", "0", "1", "200.00"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", $batsman, "1", "3", "0", "0", "33.33"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "0", "1", "0", "0", "0.00"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", $batsman, "4", "2", "1", "0", "200.00")
)

# Force every column to text formatting first so numeric-looking strings
# ("10", "200.00", ...) are stored as text, not auto-converted to numbers.
$ws.Range("A1:K5").NumberFormat = "@"

for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c).Value = $row[$c - 1]
    }
}
